$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect numeric-looking price/percentage text in columns D and E from Excel's
# automatic number conversion by temporarily marking the ranges as Text before
# writing the values, then clearing the formatting again so the cells keep the
# same (default) style as in the source file.
$colD = $ws.Range("D2:D51")
$colE = $ws.Range("E2:E51")
$colD.NumberFormat = "@"
$colE.NumberFormat = "@"

$ws.Range("D2").Value = '66.692.76'
$ws.Range("E2").Value = '  +2.72%  '
$ws.Range("D3").Value = '3.206.91'
$ws.Range("E3").Value = '  +1.64%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '602.98'
$ws.Range("E5").Value = '  +3.99%  '
$ws.Range("D6").Value = '156.42'
$ws.Range("E6").Value = '  +4.66%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").Value = '0.559'
$ws.Range("E8").Value = '  +6.33%  '
$ws.Range("D9").Value = '3.206.08'
$ws.Range("E9").Value = '  +1.65%  '
$ws.Range("E10").Value = '  +1.79%  '
$ws.Range("D11").Value = '5.84'
$ws.Range("E11").Value = '  -4.95%  '
$ws.Range("D12").Value = '0.521'
$ws.Range("E12").Value = '  +3.95%  '
$ws.Range("D13").Value = '0.0000270'
$ws.Range("E13").Value = '  +2.49%  '
$ws.Range("D14").Value = '39.29'
$ws.Range("E14").Value = '  +5.69%  '
$ws.Range("D15").Value = '3.732.42'
$ws.Range("E15").Value = '  +1.59%  '
$ws.Range("D16").Value = '66.700.22'
$ws.Range("E16").Value = '  +2.78%  '
$ws.Range("E17").Value = '  +5.16%  '
$ws.Range("D18").Value = '3.206.16'
$ws.Range("E18").Value = '  +1.52%  '
$ws.Range("D19").Value = '526.54'
$ws.Range("E19").Value = '  +4.45%  '
$ws.Range("E20").Value = '  +0.74%  '
$ws.Range("D21").Value = '15.54'
$ws.Range("E21").Value = '  +4.22%  '
$ws.Range("D22").Value = '0.743'
$ws.Range("E22").Value = '  +3.95%  '
$ws.Range("D23").Value = '8.21'
$ws.Range("E23").Value = '  +6.24%  '
$ws.Range("D24").Value = '15.06'
$ws.Range("E24").Value = '  -0.66%  '
$ws.Range("D25").Value = '85.85'
$ws.Range("E25").Value = '  +1.52%  '
$ws.Range("E26").Value = '  +0.22%  '
$ws.Range("D27").Value = '9.28'
$ws.Range("E27").Value = '  +3.12%  '
$ws.Range("E28").Value = '  +3.19%  '
$ws.Range("E29").Value = '  +9.06%  '
$ws.Range("D30").Value = '2.98'
$ws.Range("E30").Value = '  +6.71%  '
$ws.Range("D31").Value = '7.03'
$ws.Range("E31").Value = '  +10.64%  '
$ws.Range("D32").Value = '28.38'
$ws.Range("E32").Value = '  +2.84%  '
$ws.Range("D33").Value = '1.23'
$ws.Range("E33").Value = '  +2.95%  '
$ws.Range("E34").Value = '  +0.02%  '
$ws.Range("D35").Value = '6.58'
$ws.Range("E35").Value = '  +1.51%  '
$ws.Range("D36").Value = '512.07'
$ws.Range("E36").Value = '  +6.78%  '
$ws.Range("D37").Value = '54.92'
$ws.Range("E37").Value = '  +0.22%  '
$ws.Range("E38").Value = '  +1.83%  '
$ws.Range("D39").Value = '0.0428'
$ws.Range("E39").Value = '  +3.08%  '
$ws.Range("D40").Value = '0.127'
$ws.Range("E40").Value = '  +8.63%  '
$ws.Range("E41").Value = '  +2.26%  '
$ws.Range("D42").Value = '2.90'
$ws.Range("E42").Value = '  -0.57%  '
$ws.Range("D43").Value = '0.0₃0687'
$ws.Range("E43").Value = '  +16.10%  '
$ws.Range("E44").Value = '  +6.71%  '
$ws.Range("D45").Value = '2.46'
$ws.Range("E45").Value = '  +1.37%  '
$ws.Range("D46").Value = '2.896.97'
$ws.Range("E46").Value = '  -3.13%  '
$ws.Range("D47").Value = '28.62'
$ws.Range("E47").Value = '  +1.26%  '
$ws.Range("D48").Value = '2.76'
$ws.Range("E48").Value = '  +10.82%  '
$ws.Range("E49").Value = '  +3.82%  '
$ws.Range("B50").Value = 'USDe'
$ws.Range("C50").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D50").Value = '0.999'
$ws.Range("E50").Value = '  -0.01%  '
$ws.Range("B51").Value = 'ThetaToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D51").Value = '2.36'
$ws.Range("E51").Value = '  +5.10%  '

$colD.ClearFormats()
$colE.ClearFormats()

